# This workbook's last data block (rows 435-465 on Sheet1) is a weekly log of
# Frutilla (strawberry) price records for "Macroferia Regional de Talca".
# The commit adds one new weekly record. In terms of raw rows, this means:
#   - A new row is inserted right after row 435 (shifting old rows 436-465
#     down to 437-466, and the dimension grows from T465 to T466).
#   - The new row 436 receives the data that row 435 held *before* this edit.
#   - Row 435 itself is updated in place: a new date (Fecha) and a new
#     origin (Origen) value, representing the newest record, are written
#     while all of its other fields stay the same.
#
# Net effect: before[435] becomes after[436] (shifted down), and after[435]
# is before[435] with just D (Fecha) and R (Origen) updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old row 436; this pushes old rows 436-465 to
# 437-466 and grows the used range to A1:T466, matching the dimension change
# in the diff.
$ws.Rows(436).Insert()

# Populate the newly inserted row 436 with the values that row 435 held
# prior to this edit (this is exactly what the diff shows row 436 becoming).
$ws.Range("A436").Value = 5
$ws.Range("B436").Value = 'Macroferia Regional de Talca'
$ws.Range("C436").Value = 'Maule'
$ws.Range("D436").Value = 44560
$ws.Range("E436").Value = 7
$ws.Range("F436").Value = 'Fruta'
$ws.Range("G436").Value = 100101
$ws.Range("H436").Value = 'Berries'
$ws.Range("I436").Value = 100112025
$ws.Range("J436").Value = 'Frutilla'
$ws.Range("K436").Value = 'Sin especificar'
$ws.Range("L436").Value = 'Especial'
$ws.Range("M436").Value = 300
$ws.Range("N436").Value = 7000
$ws.Range("O436").Value = 7000
$ws.Range("P436").Value = 7000
$ws.Range("Q436").Value = '$/bandeja 7 kilos'
$ws.Range("R436").Value = 'Región del Maule'
$ws.Range("S436").Value = 1000
$ws.Range("T436").Value = 7

# Update row 435 in place with the newest record's date and origin.
$ws.Range("D435").Value = 44585
$ws.Range("R435").Value = 'Provincia de Melipilla'
